$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert two new rows at row 194 (existing rows 194.. shift down by 2)
$ws.Rows.Item(194).Resize(2).Insert()

# Row 194: Europe | United Kingdom | Guernsey
$ws.Cells.Item(194, 1).Value = "Europe"
$ws.Cells.Item(194, 2).Value = "United Kingdom"
$ws.Cells.Item(194, 3).Value = "Guernsey"

# Row 195: Europe | United Kingdom | Jersey
$ws.Cells.Item(195, 1).Value = "Europe"
$ws.Cells.Item(195, 2).Value = "United Kingdom"
$ws.Cells.Item(195, 3).Value = "Jersey"

# Apply the "wrap text / vertical center" style (style index 3 in styles.xml)
# to columns B and C of the two new rows, matching the source workbook.
$ws.Range("B194:C195").VerticalAlignment = -4108
$ws.Range("B194:C195").WrapText = $true

# Update the visible window / selection to match the post-edit state
$ws.Application.ActiveWindow.ScrollRow = 180
$ws.Range("A194:C195").Select()
